$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Stays"
$ws.Range("D2").Value = "after"
$ws.Range("D3").Value = "edit"

$ws.Range("D2:D3").Select()
